$wb = $excel.ActiveWorkbook

# --- clientes sheet: fill row 2 with data ---
$clientes = $wb.Worksheets.Item("clientes")
$clientes.Range("A2").Value = "activo"
$clientes.Range("B2").Value = "17234876-8"
$clientes.Range("C2").Value = "Sigo"
$clientes.Range("D2").Value = "Kosovo"
$clientes.Range("E2").Value = "Cloac"
$clientes.Range("F2").Value = "88873234"
$clientes.Range("H2").Value = "Test"

# --- ruta_actual sheet: remove trailing empty row 4 ---
$ruta = $wb.Worksheets.Item("ruta_actual")
$ruta.Rows.Item(4).Delete()
